# This script applies the "Add files via upload" commit:
#  - The sheet was regenerated with a newer extraction timestamp
#    (2024-06-07 08:51:41 -> 2024-06-10 10:29:25), which is reflected in:
#      * the worksheet tab name (IClientBalance-20240607-085141- -> IClientBalance-20240610-102925-)
#      * every "Dt. Referencia" date in column G (45450 -> 45453, i.e. 2024-06-07 -> 2024-06-10)
#  - A handful of accounts got updated balances (columns D, E and H) reflecting
#    the newer snapshot of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the newer extraction timestamp.
$ws.Name = "IClientBalance-20240610-102925-"

# Every row's "Dt. Referencia" (column G) moves from 45450 (2024-06-07) to 45453 (2024-06-10).
$ws.Range("G2:G257").Value = 45453

# Updated balances for the rows whose underlying figures changed in this snapshot.
$ws.Range("D5").Value = -10523.71
$ws.Range("E5").Value = 0
$ws.Range("H5").Value = -10523.71

$ws.Range("D8").Value = 1736.87
$ws.Range("E8").Value = 0
$ws.Range("H8").Value = 1736.87

$ws.Range("D15").Value = -12540.69
$ws.Range("E15").Value = 0
$ws.Range("H15").Value = -12540.69

$ws.Range("D17").Value = -4667.89
$ws.Range("E17").Value = 0
$ws.Range("H17").Value = -4667.89

$ws.Range("D42").Value = -5702.56
$ws.Range("E42").Value = 0
$ws.Range("H42").Value = -5702.56

$ws.Range("D48").Value = 406.89
$ws.Range("H48").Value = 406.89

$ws.Range("D51").Value = 982.12
$ws.Range("H51").Value = 982.12

$ws.Range("D54").Value = 844.29
$ws.Range("H54").Value = 844.29

$ws.Range("D56").Value = 1185.15
$ws.Range("H56").Value = 1185.15

$ws.Range("D57").Value = -351.76
$ws.Range("E57").Value = 0
$ws.Range("H57").Value = -351.76

$ws.Range("D59").Value = -8155.01
$ws.Range("E59").Value = 0
$ws.Range("H59").Value = -8155.01

$ws.Range("D98").Value = -8198.22
$ws.Range("E98").Value = 0
$ws.Range("H98").Value = -8198.22

$ws.Range("D102").Value = 466.94
$ws.Range("H102").Value = 466.94

$ws.Range("D103").Value = -23333.66
$ws.Range("E103").Value = 0
$ws.Range("H103").Value = -23333.66

$ws.Range("D106").Value = 341.42
$ws.Range("H106").Value = 341.42

$ws.Range("D107").Value = -27132.82
$ws.Range("E107").Value = 0
$ws.Range("H107").Value = -27132.82

$ws.Range("D112").Value = 499.1
$ws.Range("H112").Value = 499.1

$ws.Range("D118").Value = 106.47
$ws.Range("H118").Value = 106.47

$ws.Range("D131").Value = -3138.34
$ws.Range("E131").Value = 0
$ws.Range("H131").Value = -3138.34

$ws.Range("D141").Value = -30160.28
$ws.Range("E141").Value = 0
$ws.Range("H141").Value = -30160.28

$ws.Range("D155").Value = 250.95
$ws.Range("E155").Value = 0
$ws.Range("H155").Value = 250.95

$ws.Range("D162").Value = 127648.22
$ws.Range("H162").Value = 127648.22

$ws.Range("D168").Value = -1223.12
$ws.Range("E168").Value = 0
$ws.Range("H168").Value = -1223.12

$ws.Range("D223").Value = 336.19
$ws.Range("H223").Value = 336.19

$ws.Range("D226").Value = -7705.99
$ws.Range("E226").Value = 0
$ws.Range("H226").Value = -7705.99

$ws.Range("D240").Value = -7651.98
$ws.Range("E240").Value = 0
$ws.Range("H240").Value = -7651.98

$ws.Range("D245").Value = 78145.99
$ws.Range("H245").Value = 78145.99
